# Generate Report for Handoff
# Adds a new handed-off file (c80df167-980f-4282-a33b-6947934428ba.md) as
# row 3 on the "Overview", "zh-cn" and "de-de" worksheets, growing each
# table by one row.

$wb = $excel.ActiveWorkbook

$commitSha = "5488b565540b253f23f4553ee92b3cdbbceaaaff"
$newMd = "c80df167-980f-4282-a33b-6947934428ba.md"
$newMdPath = "e2e\" + $newMd
$newUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newMd"

$hyperlinkColor = 15570276   # OLE (BGR) packing of RGB #6495ED, matches existing HyperLink style
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Set-BlankCell($range) {
    # Forces the engine to materialize an (empty) cell at this address
    # without changing its content, mirroring the workbook's existing
    # "" shared-string entries used for not-applicable values.
    $range.Font.Bold = $false
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMd
$wsOverview.Range("B3").Value = $newMdPath
$wsOverview.Range("B3").Font.Underline = 2
$wsOverview.Range("B3").Font.Color = $hyperlinkColor
$wsOverview.Range("C3").Value = ".md"
Set-BlankCell $wsOverview.Range("D3")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-22 10:20:24"
$wsOverview.Range("G3").NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newUrl, "", "", $newMdPath)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = $newMd
$wsZhCn.Range("A3").Font.Underline = 2
$wsZhCn.Range("A3").Font.Color = $hyperlinkColor
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "c80df167-980f-4282-a33b-6947934428ba.d3c88c913e94d49e51f39205a343e3792052712f.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-22 10:20:18"
$wsZhCn.Range("H3").NumberFormat = $dateFormat
Set-BlankCell $wsZhCn.Range("I3")
Set-BlankCell $wsZhCn.Range("J3")
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dateFormat
Set-BlankCell $wsZhCn.Range("L3")
$wsZhCn.Range("M3").Value = "True"
Set-BlankCell $wsZhCn.Range("N3")
$wsZhCn.Range("O3").Value = "False"
Set-BlankCell $wsZhCn.Range("P3")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newUrl, "", "", $newMd)

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = $newMd
$wsDeDe.Range("A3").Font.Underline = 2
$wsDeDe.Range("A3").Font.Color = $hyperlinkColor
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "c80df167-980f-4282-a33b-6947934428ba.d3c88c913e94d49e51f39205a343e3792052712f.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-22 10:20:24"
$wsDeDe.Range("H3").NumberFormat = $dateFormat
Set-BlankCell $wsDeDe.Range("I3")
Set-BlankCell $wsDeDe.Range("J3")
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dateFormat
Set-BlankCell $wsDeDe.Range("L3")
$wsDeDe.Range("M3").Value = "True"
Set-BlankCell $wsDeDe.Range("N3")
$wsDeDe.Range("O3").Value = "False"
Set-BlankCell $wsDeDe.Range("P3")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newUrl, "", "", $newMd)

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

Write-Output "Handback row added for $newMd"
